$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 held "ROSHAN MONEY EXCHANGE" in column A (Nome). Correct the casing
# of the name in A5, and add the original all-caps spelling as the
# "Emitente" (issuer) entry in B5.
$ws.Range("A5").Value = "ROSHAN Money EXCHANGE"
$ws.Range("B5").Value = "ROSHAN MONEY EXCHANGE"

# End up with B5 selected, matching the saved cursor position.
$null = $ws.Range("B5").Select()
